$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Remove the "NextTransaction" workblock entries (rows 11-12: wbNextTransaction_Type /
# wbNextTransaction_SuppressSuccessful). NextTransaction.xaml was deleted and replaced
# with actions in transitions, so this workblock config row pair is no longer needed.
$ws.Rows("11:12").Delete()

$null = $ws.Range("C22").Select()
